$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.306.61"
$ws.Range("E2").Value = "  +0.54%  "

$ws.Range("D3").Value = "1.801.36"
$ws.Range("E3").Value = "  -0.20%  "

$ws.Range("E4").Value = "  +0.33%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "326.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.92%  "

$ws.Range("E6").Value = "  +0.44%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4428"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +12.00%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3721"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +8.00%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.69"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.35%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.149"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.26%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07518"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.55%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.52"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.55%  "

$ws.Range("E13").Value = "  +0.37%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.690"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.29%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.285"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.41%  "

$ws.Range("D16").Value = "1.797.24"
$ws.Range("E16").Value = "  -0.10%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001089"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.46%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06786"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.33%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "80.91"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.90%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.34%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.43"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.26%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.318"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.22%  "

$ws.Range("D23").Value = "28.298.02"
$ws.Range("E23").Value = "  +0.37%  "

$ws.Range("E24").Value = "  -0.80%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.417"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.80%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "20.41"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.59%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "152.59"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.62%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.365"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.54%  "

$ws.Range("D29").Value = "2.001.44"
$ws.Range("E29").Value = "  +0.04%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "132.69"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.01%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.234"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.17%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.018"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.05%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.822"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.91%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.09325"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.44%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.2308"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.61%  "

$ws.Range("E36").Value = "  -2.02%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06311"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.07%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02327"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.05%  "

$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6567"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.59%  "

$ws.Range("B40").Value = "InternetComputer(DFINITY)"
$ws.Range("C40").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.149"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.59%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.462"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.25%  "

$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.205"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.80%  "

$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.183"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.62%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9995"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.24%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.96"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.77%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6062"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.58%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.784"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.65%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "129.01"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.39%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.033"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.35%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07112"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.19%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.157"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.37%  "
